$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy a clean, default-styled cell once; reused to restore format after forcing text cells below.
$ws.Range("A1").Copy()

$ws.Range('D2').Value = '68.281.38'
$ws.Range('E2').Value = '  +2.02%  '
$ws.Range('D3').Value = '3.629.61'
$ws.Range('E3').Value = '  +0.92%  '
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '196.42'
$ws.Range('D5').PasteSpecial(-4122)
$ws.Range('E5').Value = '  +7.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '579.40'
$ws.Range('D6').PasteSpecial(-4122)
$ws.Range('E6').Value = '  -1.01%  '
$ws.Range('D7').Value = '3.622.99'
$ws.Range('E7').Value = '  +0.98%  '
$ws.Range('E8').Value = '  +2.74%  '
$ws.Range('E9').Value = '  -0.10%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.678'
$ws.Range('D10').PasteSpecial(-4122)
$ws.Range('E10').Value = '  +1.01%  '
$ws.Range('E11').Value = '  +7.53%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '55.70'
$ws.Range('D12').PasteSpecial(-4122)
$ws.Range('E12').Value = '  +4.26%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000298'
$ws.Range('D13').PasteSpecial(-4122)
$ws.Range('E13').Value = '  +20.39%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.13'
$ws.Range('D14').PasteSpecial(-4122)
$ws.Range('E14').Value = '  +2.66%  '
$ws.Range('D15').Value = '4.204.16'
$ws.Range('E15').Value = '  +0.61%  '
$ws.Range('D16').Value = '3.627.57'
$ws.Range('E16').Value = '  +0.73%  '
$ws.Range('E17').Value = '  +0.42%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.57'
$ws.Range('D18').PasteSpecial(-4122)
$ws.Range('E18').Value = '  +3.79%  '
$ws.Range('D19').Value = '68.170.21'
$ws.Range('E19').Value = '  +2.14%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.52'
$ws.Range('D20').PasteSpecial(-4122)
$ws.Range('E20').Value = '  +1.36%  '
$ws.Range('E21').Value = '  +1.71%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '403.81'
$ws.Range('D22').PasteSpecial(-4122)
$ws.Range('E22').Value = '  +2.96%  '
$ws.Range('E23').Value = '  -1.86%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.64'
$ws.Range('D24').PasteSpecial(-4122)
$ws.Range('E24').Value = '  +22.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '86.10'
$ws.Range('D25').PasteSpecial(-4122)
$ws.Range('E25').Value = '  +1.64%  '
$ws.Range('E26').Value = '  +4.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.62'
$ws.Range('D27').PasteSpecial(-4122)
$ws.Range('E27').Value = '  +4.14%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.89'
$ws.Range('D28').PasteSpecial(-4122)
$ws.Range('E28').Value = '  +8.99%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.14'
$ws.Range('D29').PasteSpecial(-4122)
$ws.Range('E29').Value = '  +1.57%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.15'
$ws.Range('D30').PasteSpecial(-4122)
$ws.Range('E30').Value = '  +20.57%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.16'
$ws.Range('D31').PasteSpecial(-4122)
$ws.Range('E31').Value = '  +3.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '31.74'
$ws.Range('D32').PasteSpecial(-4122)
$ws.Range('E32').Value = '  +2.46%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '693.77'
$ws.Range('D33').PasteSpecial(-4122)
$ws.Range('E33').Value = '  +20.21%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '12.22'
$ws.Range('D34').PasteSpecial(-4122)
$ws.Range('E34').Value = '  +3.58%  '
$ws.Range('E35').Value = '  +6.33%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '64.90'
$ws.Range('D36').PasteSpecial(-4122)
$ws.Range('E36').Value = '  -1.20%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '42.67'
$ws.Range('D37').PasteSpecial(-4122)
$ws.Range('E37').Value = '  +3.25%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.418'
$ws.Range('D38').PasteSpecial(-4122)
$ws.Range('E38').Value = '  +12.59%  '
$ws.Range('B39').Value = 'Dai'
$ws.Range('C39').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.00'
$ws.Range('D39').PasteSpecial(-4122)
$ws.Range('E39').Value = '  +0.14%  '
$ws.Range('B40').Value = 'PEPE'
$ws.Range('C40').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D40').Value = '0.0₃0800'
$ws.Range('E40').Value = '  +11.36%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.87'
$ws.Range('D41').PasteSpecial(-4122)
$ws.Range('E41').Value = '  +22.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.14'
$ws.Range('D42').PasteSpecial(-4122)
$ws.Range('E42').Value = '  +14.72%  '
$ws.Range('E43').Value = '  +2.69%  '
$ws.Range('D44').Value = '3.145.89'
$ws.Range('E44').Value = '  +17.73%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.998'
$ws.Range('D45').PasteSpecial(-4122)
$ws.Range('E45').Value = '  -0.22%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.92'
$ws.Range('D46').PasteSpecial(-4122)
$ws.Range('E46').Value = '  +27.55%  '
$ws.Range('E47').Value = '  +3.94%  '
$ws.Range('E48').Value = '  +1.40%  '
$ws.Range('E49').Value = '  +6.13%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.12'
$ws.Range('D50').PasteSpecial(-4122)
$ws.Range('E50').Value = '  +1.99%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '142.74'
$ws.Range('D51').PasteSpecial(-4122)
$ws.Range('E51').Value = '  +2.12%  '
